$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests")

# Pass 1: touch cells in the exact order new shared strings were first introduced
# (matches the author-observed sharedStrings.xml append order).
$ws.Range("A116").Value = 'Test: Spiel mit zwei KIs'
$ws.Range("A117").Value = 'Der Nutzer wählt für Spieler 1 "KI-Elimination", für Spieler 2 "KI-Fehlerrückführung" und klickt auf das Feld "Spiel starten".'
$ws.Range("B117").Value = 'Der Graph ist im Ausgangszustand. Es stehen zwei Beschriftungen an den Kanten. Die Beschriftungen an den Kanten sind rot und blau gefärbt.'
$ws.Range("A118").Value = 'Der Nutzer klickt auf das "NextMove"-Symbol'
$ws.Range("B119").Value = 'Zufällige freie Felder werden von den KIs besetzt.'
$ws.Range("A120").Value = 'Der Nutzer klickt auf das "Play"-Symbol, bevor das Spiel beendet ist.'
$ws.Range("B120").Value = 'Keine weiteren Felder werden besetzt.'
$ws.Range("A121").Value = 'Der Nutzer klickt auf das "Play"-Symbol und wartet bis das Spiel beendet ist. '
$ws.Range("B122").Value = 'Der Belohnungsbildschirm wird für die Evaluationsstrategie Elimination aufgerufen.'
$ws.Range("B124").Value = 'Der Belohnungsbildschirm wird für die Evaluationsstrategie Fehlerrückführung aufgerufen.'
$ws.Range("A127").Value = 'Der Nutzer wählt für Spieler 2 "KI-Elimination".'
$ws.Range("A128").Value = 'Der Nutzer klickt auf "Spiel starten".'
$ws.Range("B128").Value = 'Es wird in die Spielansicht gewechselt. Das Feld ist leer und der Graph ist im Ausgangszustand.'
$ws.Range("A129").Value = 'Der Nutzer klickt auf das "Play"-Symbol und wartet bis das Spiel beendet ist.'
$ws.Range("A131").Value = 'Der Spieler drückt auf "Weiter".'
$ws.Range("A133").Value = 'Test: KI erzeugen'
$ws.Range("A134").Value = 'Der Nutzer klickt auf das DropDown-Menü für die Auswahl von "Wähle einen KI-Typ".'
$ws.Range("B134").Value = 'Es wird im DropDown-Menü eine Liste aller möglichen Spieler angezeigt. Die Optionen sind: "Elimination", "Fehlerrückführung", "Zufällig"'
$ws.Range("A135").Value = 'Der Nutzer wählt die Option "Elimination" im DropDown-Menü aus.'
$ws.Range("B135").Value = 'Im Feld des DropDown-Menüs für "Wähle einen KI-Typ" steht "Elimination".'
$ws.Range("A136").Value = 'Der Nutzer klickt auf das Feld "Wähle einen Namen für die KI".'
$ws.Range("B136").Value = 'Das Feld ist beschreibbar.'
$ws.Range("A137").Value = 'Der Nutzer schreibt "KI 1" in das Feld.'
$ws.Range("B137").Value = 'Im Feld "Wähle einen Namen für die KI" steht "KI 1".'
$ws.Range("A138").Value = 'Der Nutzer klickt auf "Erstelle eine neue KI".'
$ws.Range("B138").Value = 'In der Liste der KIs wird die neue KI angezeigt.'
$ws.Range("A140").Value = 'Der Nutzer wählt die Option "Fehlerrückführung" im DropDown-Menü aus.'
$ws.Range("B140").Value = 'Im Feld des DropDown-Menüs für "Wähle einen KI-Typ" steht "Fehlerrückführung".'
$ws.Range("A142").Value = 'Der Nutzer schreibt "KI 2" in das Feld.'
$ws.Range("B142").Value = 'Im Feld "Wähle einen Namen für die KI" steht "KI 2".'
$ws.Range("B143").Value = 'In der Liste wird die neue KI angezeigt.'
$ws.Range("A144").Value = 'Der Nutzer wählt für Spieler 1 "KI 1", für Spieler 2 "KI 2".'
$ws.Range("B144").Value = 'In den beiden Feldern steht "KI 1" und "KI 2".'
$ws.Range("A145").Value = 'Ein Spiel wird komplett durchgespielt und die KIs werden belohnt.'
$ws.Range("A146").Value = 'Der Nutzer klickt im Startbildschirm auf den Knopf "Zurücksetzen" neben der "KI 2".'
$ws.Range("A147").Value = 'Der Nutzer startet ein neues Spiel.'
$ws.Range("B147").Value = 'Die "KI 2" hat die gleichen Gewichte, wie beim ersten durchspielen.'
$ws.Range("A149").Value = 'Test: Slider'
$ws.Range("A150").Value = 'Der Nutzer wählt für Spieler 1 "KI-Elimination", für Spieler 2 "KI-Elimination" und klickt auf das Feld "Spiel starten".'
$ws.Range("B150").Value = 'In der Spielansicht existiert ein Slider, welcher auf ganz links eingestellt ist.'
$ws.Range("A153").Value = 'Der Nutzer bewegt den Slider nach ganz rechts.'
$ws.Range("B151").Value = 'Der Slider bewegt sich mit der Maus.'
$ws.Range("A151").Value = 'Der Nutzer klickt auf den Slider, hält die Maus gedrückt und zieht sie nach ganz rechts und wieder nach ganz links.'
$ws.Range("B153").Value = 'Die Geschwindigkeit, mit der die KI die Züge macht erhöht sich.'
$ws.Range("A154").Value = 'Der Nutzer wartet, bis das Spiel beendet ist, klickt auf "Weiter", "Belohnung ausführen", "Weiter" und auf "Spiel starten".'
$ws.Range("B154").Value = 'Der Slider ist auf ganz rechts eingestellt und die KI macht ihre Züge weiterhin sehr schnell.'
$ws.Range("A156").Value = 'Test: Speedrun'
$ws.Range("A158").Value = 'Der Nutzer stellt den Slider auf ganz rechts.'
$ws.Range("B161").Value = 'Die KIs machen sehr schnell ihre Züge und es werden mehrere Spiele hintereinander ausgeführt.'

# Pass 2: remaining cells that reuse already-known shared strings
$ws.Range("A119").Value = 'Der Nutzer klickt auf das "Play"-Symbol.'
$ws.Range("A122").Value = 'Der Nutzer klickt auf "Weiter".'
$ws.Range("A123").Value = 'Der Nutzer klickt auf das Feld "Belohnung ausführen".'
$ws.Range("A124").Value = 'Der Nutzer klickt auf "Weiter".'
$ws.Range("A125").Value = 'Der Nutzer klickt auf das Feld "Belohnung ausführen".'
$ws.Range("A126").Value = 'Der Nutzer klickt auf "Weiter".'
$ws.Range("B126").Value = 'Es wird in die Startansicht gewechselt.'
$ws.Range("B129").Value = 'Der Belohnungsbildschirm wird für die Evaluationsstrategie Elimination aufgerufen.'
$ws.Range("A130").Value = 'Der Nutzer klickt auf das Feld "Belohnung ausführen".'
$ws.Range("B131").Value = 'Es wird in die Startansicht gewechselt.'
$ws.Range("A139").Value = 'Der Nutzer klickt auf das DropDown-Menü für die Auswahl von "Wähle einen KI-Typ".'
$ws.Range("B139").Value = 'Es wird im DropDown-Menü eine Liste aller möglichen Spieler angezeigt.'
$ws.Range("A141").Value = 'Der Nutzer klickt auf das Feld "Wähle einen Namen für die KI".'
$ws.Range("B141").Value = 'Das Feld ist beschreibbar.'
$ws.Range("A143").Value = 'Der Nutzer klickt auf "Erstelle eine neue KI".'
$ws.Range("A152").Value = 'Der Nutzer klickt auf das "Play"-Symbol.'
$ws.Range("A157").Value = 'Der Nutzer wählt für Spieler 1 "KI-Elimination", für Spieler 2 "KI-Elimination" und klickt auf das Feld "Spiel starten".'
$ws.Range("A159").Value = 'Der Nutzer klickt auf die Checkbox "Startansicht überspringen".'
$ws.Range("A160").Value = ' Der Nutzer klickt auf die Checkbox "Belohnungsansicht überspringen".'
$ws.Range("A161").Value = 'Der Nutzer klickt auf das "Play"-Symbol.'

# Bold section headers (matches existing "Test: ..." header style)
$ws.Range("A116").Font.Bold = $true
$ws.Range("A133").Font.Bold = $true
$ws.Range("A149").Font.Bold = $true
$ws.Range("A156").Font.Bold = $true

$ws.Range("A162").Select()
$excel.ActiveWindow.ScrollRow = 139
$excel.ActiveWindow.ScrollColumn = 1
